$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number, week-covering dates) ---
$ws.Range("A8").Value = "Volume 30   Number  49"
$ws.Range("C9").Value = "Report Covering the Week  12/4/2023  Through  12/10/2023"

# --- Cells whose type flips between number and text ("0" / "***.*") ---
# Row 14 (Murder) is never touched by this edit, so its cells are safe formatting
# templates: C14 = style 14 / text "0", E14 = style 14 / text "***.*", I14 = style 15 / number.
# Copying reproduces the exact target style; for the text targets the copied value is
# already correct ("0"/"***.*"), for the numeric targets we overwrite the value afterwards.
$ws.Range("I14").Copy($ws.Range("C15"))
$ws.Range("C14").Copy($ws.Range("D15"))
$ws.Range("E14").Copy($ws.Range("E15"))
$ws.Range("I14").Copy($ws.Range("F15"))
$ws.Range("I14").Copy($ws.Range("C18"))
$ws.Range("C14").Copy($ws.Range("D18"))
$ws.Range("E14").Copy($ws.Range("E18"))
$ws.Range("I14").Copy($ws.Range("C20"))
$ws.Range("I14").Copy($ws.Range("C26"))
$ws.Range("C14").Copy($ws.Range("D26"))
$ws.Range("E14").Copy($ws.Range("E26"))
$ws.Range("I14").Copy($ws.Range("F26"))
$ws.Range("C14").Copy($ws.Range("D27"))
$ws.Range("E14").Copy($ws.Range("E27"))

# --- Set final values for cells that keep/become numeric, and the remaining plain value changes ---
$ws.Range("C15").Value = 1
$ws.Range("F15").Value = 1
$ws.Range("H15").Value = -66.666666666666
$ws.Range("I15").Value = 10
$ws.Range("K15").Value = -37.5
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -37.5
$ws.Range("N15").Value = -72.222222222222
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 200
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = -10
$ws.Range("I16").Value = 159
$ws.Range("J16").Value = 170
$ws.Range("K16").Value = -6.470588235294
$ws.Range("L16").Value = 54.368932038834
$ws.Range("M16").Value = -38.132295719844
$ws.Range("N16").Value = -84.334975369458
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = 14.285714285714
$ws.Range("F17").Value = 26
$ws.Range("G17").Value = 26
$ws.Range("I17").Value = 352
$ws.Range("J17").Value = 359
$ws.Range("K17").Value = -1.949860724233
$ws.Range("L17").Value = 6.344410876132
$ws.Range("M17").Value = 72.549019607843
$ws.Range("N17").Value = -50.282485875706
$ws.Range("C18").Value = 1
$ws.Range("I18").Value = 100
$ws.Range("K18").Value = -35.897435897435
$ws.Range("L18").Value = -39.393939393939
$ws.Range("M18").Value = -35.483870967741
$ws.Range("N18").Value = -89.327641408751
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = -8.333333333333
$ws.Range("F19").Value = 32
$ws.Range("G19").Value = 40
$ws.Range("H19").Value = -20
$ws.Range("I19").Value = 436
$ws.Range("J19").Value = 500
$ws.Range("K19").Value = -12.8
$ws.Range("L19").Value = 8.728179551122
$ws.Range("M19").Value = -17.424242424242
$ws.Range("N19").Value = -33.738601823708
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = -83.333333333333
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = -45.454545454545
$ws.Range("I20").Value = 105
$ws.Range("J20").Value = 116
$ws.Range("K20").Value = -9.482758620689
$ws.Range("L20").Value = 52.173913043478
$ws.Range("M20").Value = -2.777777777777
$ws.Range("N20").Value = -90.59982094897
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = -3.846153846153
$ws.Range("F21").Value = 79
$ws.Range("G21").Value = 103
$ws.Range("H21").Value = -23.300970873786
$ws.Range("I21").Value = 1167
$ws.Range("J21").Value = 1322
$ws.Range("K21").Value = -11.724659606656
$ws.Range("L21").Value = 7.162534435261
$ws.Range("M21").Value = -8.542319749216
$ws.Range("N21").Value = -74.014696058784
$ws.Range("F22").Value = 1
$ws.Range("H22").Value = -50
$ws.Range("L22").Value = 28.571428571428
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 13
$ws.Range("G23").Value = 19
$ws.Range("H23").Value = -31.578947368421
$ws.Range("I23").Value = 159
$ws.Range("J23").Value = 172
$ws.Range("K23").Value = -7.558139534883
$ws.Range("L23").Value = -8.092485549132
$ws.Range("M23").Value = 38.260869565217
$ws.Range("C24").Value = 26
$ws.Range("D24").Value = 29
$ws.Range("E24").Value = -10.344827586206
$ws.Range("F24").Value = 73
$ws.Range("G24").Value = 85
$ws.Range("H24").Value = -14.117647058823
$ws.Range("I24").Value = 943
$ws.Range("J24").Value = 1157
$ws.Range("K24").Value = -18.496110630942
$ws.Range("L24").Value = 13.477737665463
$ws.Range("M24").Value = -16.989436619718
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = -20
$ws.Range("F25").Value = 41
$ws.Range("G25").Value = 36
$ws.Range("H25").Value = 13.888888888888
$ws.Range("I25").Value = 579
$ws.Range("J25").Value = 529
$ws.Range("K25").Value = 9.451795841209
$ws.Range("L25").Value = 31.590909090909
$ws.Range("M25").Value = 16.265060240963
$ws.Range("C26").Value = 1
$ws.Range("F26").Value = 1
$ws.Range("H26").Value = -66.666666666666
$ws.Range("I26").Value = 17
$ws.Range("K26").Value = -39.285714285714
$ws.Range("L26").Value = -19.047619047619
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = -50
$ws.Range("L27").Value = 42.857142857142
$ws.Range("C28").Value = 1
$ws.Range("I28").Value = 11
$ws.Range("K28").Value = -57.692307692307
$ws.Range("L28").Value = -31.25
$ws.Range("M28").Value = -50
$ws.Range("N28").Value = -83.823529411764
$ws.Range("I29").Value = 7
$ws.Range("K29").Value = -58.823529411764
$ws.Range("L29").Value = -50
$ws.Range("M29").Value = -58.823529411764
$ws.Range("N29").Value = -86.792452830188
